$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 12
$ws.Range("P2").Value = 1.33
$ws.Range("Q2").Value = 3.25
$ws.Range("V2").Value = 23
$ws.Range("Z2").Value = 12
$ws.Range("AB2").Value = 19
$ws.Range("AD2").Value = 301
$ws.Range("AE2").Value = 8
$ws.Range("AF2").Value = 8
$ws.Range("AJ2").Value = 29
# Row 3
$ws.Range("H3").Value = 3
$ws.Range("P3").Value = 1.57
$ws.Range("Q3").Value = 2.25
$ws.Range("X3").Value = 29
$ws.Range("Z3").Value = 6.5
$ws.Range("AA3").Value = 6
$ws.Range("AF3").Value = 11
# Row 4
$ws.Range("N4").Value = 2.38
$ws.Range("O4").Value = 1.57
# Row 5
$ws.Range("G5").Value = 3.25
$ws.Range("I5").Value = 2.6
$ws.Range("J5").Value = 1.2
$ws.Range("K5").Value = 4.33
$ws.Range("AF5").Value = 10
# Row 7
$ws.Range("N7").Value = 2.6
$ws.Range("O7").Value = 1.48
# Row 16
$ws.Range("H16").Value = 4.85
$ws.Range("I16").Value = 9.75
$ws.Range("L16").Value = 1.26
$ws.Range("M16").Value = 3.15
$ws.Range("N16").Value = 1.78
$ws.Range("O16").Value = 1.83
$ws.Range("R16").Value = 2.35
$ws.Range("S16").Value = 1.47
$ws.Range("T16").Value = 5.7
$ws.Range("U16").Value = 5.3
$ws.Range("W16").Value = 7
$ws.Range("X16").Value = 12.5
$ws.Range("Y16").Value = 40
$ws.Range("Z16").Value = 10.25
$ws.Range("AB16").Value = 32
$ws.Range("AC16").Value = 200
$ws.Range("AE16").Value = 21
$ws.Range("AF16").Value = 70
$ws.Range("AG16").Value = 35
$ws.Range("AH16").Value = 350
$ws.Range("AI16").Value = 175
$ws.Range("AJ16").Value = 150
# Row 17
$ws.Range("G17").Value = 4.3
$ws.Range("L17").Value = 1.32
$ws.Range("M17").Value = 2.82
$ws.Range("N17").Value = 1.93
$ws.Range("O17").Value = 1.7
$ws.Range("R17").Value = 1.83
$ws.Range("S17").Value = 1.78
$ws.Range("Z17").Value = 8.5
$ws.Range("AB17").Value = 16
$ws.Range("AC17").Value = 80
$ws.Range("AD17").Value = 700
$ws.Range("AE17").Value = 6.4
$ws.Range("AF17").Value = 8
$ws.Range("AH17").Value = 14.5
$ws.Range("AI17").Value = 14.5
$ws.Range("AJ17").Value = 29
# Row 19
$ws.Range("G19").Value = 2.55
$ws.Range("I19").Value = 2.8
$ws.Range("N19").Value = 2
$ws.Range("O19").Value = 1.8
$ws.Range("AC19").Value = 41
$ws.Range("AE19").Value = 9
$ws.Range("AF19").Value = 13
$ws.Range("AJ19").Value = 29
# Row 24
$ws.Range("J24").Value = 1.12
$ws.Range("K24").Value = 5.4
$ws.Range("M24").Value = 2.32
$ws.Range("N24").Value = 2.57
$ws.Range("Q24").Value = 2.2
$ws.Range("R24").Value = 2.32
$ws.Range("T24").Value = 5
$ws.Range("U24").Value = 7.2
$ws.Range("V24").Value = 9.5
$ws.Range("Y24").Value = 50
$ws.Range("Z24").Value = 5.4
# Row 25
$ws.Range("I25").Value = 7.4
$ws.Range("J25").Value = 1.08
$ws.Range("K25").Value = 6.6
$ws.Range("L25").Value = 1.37
$ws.Range("M25").Value = 2.85
$ws.Range("N25").Value = 2.1
$ws.Range("O25").Value = 1.65
$ws.Range("Q25").Value = 2.5
$ws.Range("R25").Value = 2.37
$ws.Range("S25").Value = 1.52
$ws.Range("T25").Value = 5.1
$ws.Range("U25").Value = 5.5
$ws.Range("W25").Value = 8.5
$ws.Range("X25").Value = 14
$ws.Range("Z25").Value = 6.6
$ws.Range("AB25").Value = 27
$ws.Range("AH25").Value = 250
# Row 27
$ws.Range("G27").Value = 5.25
$ws.Range("H27").Value = 3.9
$ws.Range("I27").Value = 1.65
$ws.Range("J27").Value = 1.06
$ws.Range("K27").Value = 10
$ws.Range("N27").Value = 2
$ws.Range("O27").Value = 1.85
$ws.Range("T27").Value = 13
$ws.Range("AI27").Value = 13
# Row 29
$ws.Range("J29").Value = 1.06
$ws.Range("K29").Value = 10
$ws.Range("L29").Value = 1.3
$ws.Range("M29").Value = 3.4
$ws.Range("N29").Value = 2.05
$ws.Range("O29").Value = 1.75
# Row 30
$ws.Range("K30").Value = 13
# Row 31
$ws.Range("G31").Value = 1.55
$ws.Range("H31").Value = 3.7
$ws.Range("I31").Value = 5.75
$ws.Range("K31").Value = 9
$ws.Range("N31").Value = 2.05
$ws.Range("O31").Value = 1.75
$ws.Range("Y31").Value = 29
$ws.Range("Z31").Value = 9
$ws.Range("AA31").Value = 7.5
$ws.Range("AB31").Value = 21
$ws.Range("AI31").Value = 51
# Row 32
$ws.Range("G32").Value = 1.91
$ws.Range("H32").Value = 3.2
$ws.Range("J32").Value = 1.06
$ws.Range("K32").Value = 10
$ws.Range("L32").Value = 1.3
$ws.Range("M32").Value = 3.4
$ws.Range("N32").Value = 2.05
$ws.Range("O32").Value = 1.8
$ws.Range("P32").Value = 1.44
$ws.Range("Q32").Value = 2.63
$ws.Range("R32").Value = 1.83
$ws.Range("S32").Value = 1.83
$ws.Range("T32").Value = 7
$ws.Range("U32").Value = 9
$ws.Range("X32").Value = 17
$ws.Range("Y32").Value = 29
$ws.Range("Z32").Value = 9
$ws.Range("AD32").Value = 251
$ws.Range("AE32").Value = 11
$ws.Range("AG32").Value = 15
# Row 33
$ws.Range("L33").Value = 1.29
$ws.Range("N33").Value = 1.87
$ws.Range("P33").Value = 1.42
$ws.Range("Q33").Value = 2.65
$ws.Range("R33").Value = 1.7
$ws.Range("T33").Value = 7.4
$ws.Range("AB33").Value = 13
$ws.Range("AC33").Value = 55
$ws.Range("AE33").Value = 11.25
$ws.Range("AF33").Value = 22
# Row 34
$ws.Range("G34").Value = 6.4
$ws.Range("H34").Value = 3.9
$ws.Range("I34").Value = 1.45
$ws.Range("R34").Value = 2.27
$ws.Range("S34").Value = 1.57
$ws.Range("T34").Value = 14
$ws.Range("U34").Value = 40
$ws.Range("V34").Value = 23
$ws.Range("W34").Value = 150
$ws.Range("X34").Value = 90
$ws.Range("Y34").Value = 100
$ws.Range("AA34").Value = 8
$ws.Range("AB34").Value = 25
$ws.Range("AC34").Value = 150
$ws.Range("AE34").Value = 5.4
$ws.Range("AF34").Value = 5.8
$ws.Range("AH34").Value = 9.25
# Row 36
$ws.Range("H36").Value = 3.95
$ws.Range("I36").Value = 5.5
$ws.Range("K36").Value = 7.7
$ws.Range("L36").Value = 1.28
$ws.Range("N36").Value = 1.82
$ws.Range("O36").Value = 1.88
$ws.Range("U36").Value = 6.8
$ws.Range("Z36").Value = 7.7
$ws.Range("AB36").Value = 19
# Row 38
$ws.Range("G38").Value = 1.7
$ws.Range("I38").Value = 4.5
$ws.Range("K38").Value = 12
$ws.Range("U38").Value = 8.5
$ws.Range("AB38").Value = 15
$ws.Range("AI38").Value = 34
# Row 39
$ws.Range("N39").Value = 1.85
$ws.Range("O39").Value = 2
# Row 41
$ws.Range("G41").Value = 2.12
$ws.Range("H41").Value = 2.92
$ws.Range("I41").Value = 3.6
$ws.Range("J41").Value = 1.13
$ws.Range("K41").Value = 4.5
$ws.Range("L41").Value = 1.57
$ws.Range("M41").Value = 2.1
$ws.Range("N41").Value = 2.65
$ws.Range("O41").Value = 1.37
$ws.Range("P41").Value = 1.6
$ws.Range("Q41").Value = 2.07
$ws.Range("R41").Value = 2.25
$ws.Range("S41").Value = 1.5
$ws.Range("T41").Value = 5.1
$ws.Range("U41").Value = 8.25
$ws.Range("V41").Value = 10
$ws.Range("W41").Value = 20
$ws.Range("X41").Value = 24
$ws.Range("Y41").Value = 55
$ws.Range("Z41").Value = 4.85
$ws.Range("AA41").Value = 6.1
$ws.Range("AB41").Value = 22
$ws.Range("AC41").Value = 175
$ws.Range("AE41").Value = 7.2
$ws.Range("AF41").Value = 16.5
$ws.Range("AG41").Value = 14
$ws.Range("AH41").Value = 55
$ws.Range("AI41").Value = 45
$ws.Range("AJ41").Value = 75
# Row 42
$ws.Range("G42").Value = 2.25
$ws.Range("H42").Value = 2.92
$ws.Range("I42").Value = 3.3
$ws.Range("J42").Value = 1.1
$ws.Range("K42").Value = 6.2
$ws.Range("L42").Value = 1.47
$ws.Range("M42").Value = 2.32
$ws.Range("N42").Value = 2.37
$ws.Range("O42").Value = 1.45
$ws.Range("P42").Value = 1.52
$ws.Range("Q42").Value = 2.22
$ws.Range("R42").Value = 2.02
$ws.Range("S42").Value = 1.62
$ws.Range("T42").Value = 5.8
$ws.Range("U42").Value = 9.5
$ws.Range("V42").Value = 9.5
$ws.Range("W42").Value = 22
$ws.Range("X42").Value = 23
$ws.Range("Y42").Value = 45
$ws.Range("Z42").Value = 6.5
$ws.Range("AA42").Value = 5.8
$ws.Range("AB42").Value = 18
$ws.Range("AC42").Value = 110
$ws.Range("AE42").Value = 7.6
$ws.Range("AF42").Value = 15.5
$ws.Range("AG42").Value = 12
$ws.Range("AH42").Value = 45
$ws.Range("AI42").Value = 37
$ws.Range("AJ42").Value = 55
# Row 43
$ws.Range("G43").Value = 3.7
$ws.Range("H43").Value = 3.4
$ws.Range("I43").Value = 1.95
$ws.Range("K43").Value = 9.5
$ws.Range("AF43").Value = 9
